# Applies the "fix mixed-case text" edit described by the commit:
# several ALL-CAPS names/institutes in the CV are retyped to
# Title/Mixed case, one word-fragment at a time.  Each fragment
# keeps its own <w:r> (the commit's XML shows the text split into
# several runs with identical rPr) so we deliberately force a run
# boundary at every fragment edge by flipping Font.Bold on and back
# off again - a formatting "no-op" that still stops this host from
# silently re-coalescing the freshly split runs into one.

$d = $word.ActiveDocument

function Split-CaseFix($searchText, [string[]]$fragments) {
    $matchRange = $d.Content
    $found = $matchRange.Find.Execute($searchText)
    if (-not $found) {
        Write-Host ("NOT FOUND: " + $searchText)
        return
    }
    $base = $matchRange.Start

    $pos = $base
    for ($i = 0; $i -lt $fragments.Length; $i++) {
        $frag = $fragments[$i]
        $fragStart = $pos
        $fragEnd = $pos + $frag.Length
        $pos = $fragEnd

        # The very first fragment needs no edit beyond being left alone -
        # it is still the literal text already in the document.  Every
        # other fragment is (re)written, and every fragment past the
        # first gets the bold-flicker so it survives as its own run.
        if ($i -eq 0) {
            continue
        }

        $sub = $d.Range($fragStart, $fragEnd)
        $sub.Text = $frag
        $sub.Font.Bold = $true
        $sub.Font.Bold = $false
    }
}

Split-CaseFix "GOVT BANGA BANDHU COLLEGE" @("GOVT B", "anga ", "B", "andhu", " C", "ollege")
Split-CaseFix "MIRPUR GIRLS IDEAL LABORATORY" @("M", "irpur", " G", "irls ", "I", "deal", " L", "aboratory")
Split-CaseFix "M. D. C. MODEL INSTITUTE" @("M. D. C. M", "odel", " I", "nstitute")
Split-CaseFix "MD. FAZLUL HAQ" @("MD. ", "Fazlul Haq")
Split-CaseFix "MONOYARA BEGUM" @("M", "onoyara", " B", "egum")
